$d = $word.ActiveDocument

$pairs = @(
    @("810÷2=405, 0", "939÷8=117, 3"),
    @("748÷4=187, 0", "398÷2=199, 0"),
    @("145÷9=16, 1", "957÷2=478, 1"),
    @("831÷7=118, 5", "769÷3=256, 1"),
    @("417÷5=83, 2", "511÷6=85, 1"),
    @("252÷7=36, 0", "685÷4=171, 1"),
    @("352÷5=70, 2", "453÷9=50, 3"),
    @("787÷3=262, 1", "572÷9=63, 5"),
    @("428÷7=61, 1", "648÷3=216, 0"),
    @("615÷2=307, 1", "612÷9=68, 0"),
    @("991÷6=165, 1", "983÷5=196, 3"),
    @("568÷4=142, 0", "879÷7=125, 4"),
    @("663÷2=331, 1", "511÷5=102, 1"),
    @("830÷9=92, 2", "913÷4=228, 1"),
    @("955÷3=318, 1", "863÷8=107, 7"),
    @("533÷7=76, 1", "112÷6=18, 4"),
    @("655÷4=163, 3", "313÷8=39, 1"),
    @("459÷6=76, 3", "914÷4=228, 2"),
    @("883÷7=126, 1", "295÷5=59, 0"),
    @("563÷7=80, 3", "703÷3=234, 1"),
    @("837÷6=139, 3", "865÷6=144, 1"),
    @("120÷5=24, 0", "786÷6=131, 0"),
    @("733÷5=146, 3", "296÷6=49, 2"),
    @("868÷5=173, 3", "347÷5=69, 2"),
    @("561÷5=112, 1", "589÷8=73, 5")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
